$wb = $excel.ActiveWorkbook

# This reverts the "wrong-xlsform-col" fix on the "survey" sheet: the
# third header cell (C1) goes back from "label" to "message", and the
# active selection moves from C2 down to A3.
$ws = $wb.Worksheets.Item("survey")
$ws.Range("C1").Value = "message"
$ws.Range("A3").Select()
